$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "2021" year column (O) mirroring the existing columns,
# copying formatting from column N (the previous last year column) and
# filling in the new year's data.
$ws.Range("N4:N14").Copy() | Out-Null
$ws.Range("O4:O14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("O4").Value = 2021
$ws.Range("O5").Value = 2
$ws.Range("O6").Value = "-"
$ws.Range("O7").Value = 1
$ws.Range("O8").Value = "-"
$ws.Range("O9").Value = "-"
$ws.Range("O10").Value = "-"
$ws.Range("O11").Value = "-"
$ws.Range("O12").Value = 1
$ws.Range("O13").Value = "-"
$ws.Range("O14").Value = "-"

$ws.Range("P1").Select() | Out-Null
